$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 607.8
$ws.Range("I33").Value = 249.5
$ws.Range("K33").Value = 249.5
$ws.Range("M33").Value = -20.5
$ws.Range("H40").Value = 4432.5835
$ws.Range("I40").Value = 4228.5713
$ws.Range("K40").Value = 4228.5713
$ws.Range("M40").Value = -4053.5713
$ws.Range("H64").Value = 67082.75
$ws.Range("J64").Value = 22598.6
$ws.Range("L64").Value = 22598.6
$ws.Range("N64").Value = -23094.6
$ws.Range("H67").Value = 67082.75
$ws.Range("J67").Value = 22598.6
$ws.Range("L67").Value = 22598.6
$ws.Range("N67").Value = -24314.6
$ws.Range("H70").Value = 14287141
$ws.Range("J70").Value = 1634.3636
$ws.Range("L70").Value = 4903.0908
$ws.Range("N70").Value = -5443.0908
$ws.Range("H73").Value = 14287141
$ws.Range("J73").Value = 1634.3636
$ws.Range("L73").Value = 4903.0908
$ws.Range("N73").Value = -6775.0908
$ws.Range("H92").Value = 978
$ws.Range("I92").Value = 1027.4166
$ws.Range("K92").Value = 1027.4166
$ws.Range("M92").Value = 220.5834
$ws.Range("H132").Value = 2567862.2
$ws.Range("I132").Value = 3653.853
$ws.Range("K132").Value = 10961.559
$ws.Range("M132").Value = -8431.559000000001
$ws.Range("H137").Value = 11130.652
$ws.Range("I137").Value = 19137.25
$ws.Range("K137").Value = 57411.75
$ws.Range("M137").Value = -54861.75
$ws.Range("H138").Value = 323625.62
$ws.Range("I138").Value = 599273.2
$ws.Range("K138").Value = 1797819.6
$ws.Range("M138").Value = -1792679.6
$ws.Range("H141").Value = 6613.7646
$ws.Range("I141").Value = 6428.933
$ws.Range("K141").Value = 19286.799
$ws.Range("M141").Value = -14106.799

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 128303
$ws.Range("I2").Value = 747
$ws.Range("J2").Value = 170821.67
$ws.Range("K2").Value = 747
$ws.Range("L2").Value = 170821.67
$ws.Range("M2").Value = -634
$ws.Range("N2").Value = -171047.67
$ws.Range("H32").Value = 6950.276
$ws.Range("I32").Value = 6887.0386
$ws.Range("J32").Value = 7498.3335
$ws.Range("K32").Value = 6887.0386
$ws.Range("L32").Value = 7498.3335
$ws.Range("M32").Value = -6600.0386
$ws.Range("N32").Value = -8072.3335
$ws.Range("H61").Value = 8029.3228
$ws.Range("I61").Value = 8866.218000000001
$ws.Range("K61").Value = 8866.218000000001
$ws.Range("M61").Value = -8654.218000000001
$ws.Range("H88").Value = 50001056
$ws.Range("I88").Value = 507.125
$ws.Range("J88").Value = 83334750
$ws.Range("K88").Value = 507.125
$ws.Range("L88").Value = 83334750
$ws.Range("M88").Value = -101.125
$ws.Range("N88").Value = -83335562
$ws.Range("H91").Value = 50001056
$ws.Range("I91").Value = 507.125
$ws.Range("J91").Value = 83334750
$ws.Range("K91").Value = 507.125
$ws.Range("L91").Value = 83334750
$ws.Range("M91").Value = 896.875
$ws.Range("N91").Value = -83337558
$ws.Range("H102").Value = 9506
$ws.Range("I102").Value = 13063.389
$ws.Range("K102").Value = 13063.389
$ws.Range("M102").Value = -11441.389
$ws.Range("H116").Value = 128303
$ws.Range("I116").Value = 747
$ws.Range("J116").Value = 170821.67
$ws.Range("K116").Value = 747
$ws.Range("L116").Value = 170821.67
$ws.Range("M116").Value = 1547
$ws.Range("N116").Value = -175409.67
$ws.Range("H122").Value = 1039115.44
$ws.Range("I122").Value = 4000.45
$ws.Range("K122").Value = 12001.35
$ws.Range("M122").Value = -9551.349999999999
$ws.Range("H132").Value = 3091.525
$ws.Range("I132").Value = 2871.111
$ws.Range("K132").Value = 8613.332999999999
$ws.Range("M132").Value = -6083.332999999999
$ws.Range("H136").Value = 8029.3228
$ws.Range("I136").Value = 8866.218000000001
$ws.Range("K136").Value = 26598.654
$ws.Range("M136").Value = -24048.654

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 128303
$ws.Range("I3").Value = 747
$ws.Range("J3").Value = 170821.67
$ws.Range("K3").Value = 747
$ws.Range("L3").Value = 170821.67
$ws.Range("M3").Value = -633
$ws.Range("N3").Value = -171049.67
$ws.Range("H20").Value = 3471.842
$ws.Range("I20").Value = 2044.2142
$ws.Range("K20").Value = 2044.2142
$ws.Range("M20").Value = -1797.2142
$ws.Range("H94").Value = 9411.593999999999
$ws.Range("I94").Value = 13521.85
$ws.Range("K94").Value = 13521.85
$ws.Range("M94").Value = -13070.85
$ws.Range("H105").Value = 75955.78999999999
$ws.Range("I105").Value = 113042.336
$ws.Range("K105").Value = 113042.336
$ws.Range("M105").Value = -111295.336
$ws.Range("H134").Value = 8037.5
$ws.Range("I134").Value = 8197.368
$ws.Range("K134").Value = 24592.104
$ws.Range("M134").Value = -22057.104

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6695.794
$ws.Range("I31").Value = 7709.3335
$ws.Range("J31").Value = 5058.5386
$ws.Range("K31").Value = 7709.3335
$ws.Range("L31").Value = 5058.5386
$ws.Range("M31").Value = -7414.3335
$ws.Range("N31").Value = -5648.5386
$ws.Range("H34").Value = 6695.794
$ws.Range("I34").Value = 7709.3335
$ws.Range("J34").Value = 5058.5386
$ws.Range("K34").Value = 7709.3335
$ws.Range("L34").Value = 5058.5386
$ws.Range("M34").Value = -7507.3335
$ws.Range("N34").Value = -5462.5386
$ws.Range("H58").Value = 2257.8333
$ws.Range("I58").Value = 2205.2917
$ws.Range("J58").Value = 2468
$ws.Range("K58").Value = 2205.2917
$ws.Range("L58").Value = 2468
$ws.Range("M58").Value = -2002.2917
$ws.Range("N58").Value = -2874
$ws.Range("H105").Value = 213941.4
$ws.Range("I105").Value = 266365
$ws.Range("J105").Value = 4247
$ws.Range("K105").Value = 266365
$ws.Range("L105").Value = 4247
$ws.Range("M105").Value = -264618
$ws.Range("N105").Value = -7741
$ws.Range("H122").Value = 10036.261
$ws.Range("J122").Value = 13005.571
$ws.Range("L122").Value = 39016.713
$ws.Range("N122").Value = -43916.713
$ws.Range("H132").Value = 2684
$ws.Range("I132").Value = 2552.4546
$ws.Range("K132").Value = 7657.3638
$ws.Range("M132").Value = -5127.3638
$ws.Range("H134").Value = 7343.391
$ws.Range("I134").Value = 9860
$ws.Range("J134").Value = 2624.75
$ws.Range("K134").Value = 29580
$ws.Range("L134").Value = 7874.25
$ws.Range("M134").Value = -27045
$ws.Range("N134").Value = -12944.25
$ws.Range("H136").Value = 2257.8333
$ws.Range("I136").Value = 2205.2917
$ws.Range("J136").Value = 2468
$ws.Range("K136").Value = 6615.875100000001
$ws.Range("L136").Value = 7404
$ws.Range("M136").Value = -4065.875100000001
$ws.Range("N136").Value = -12504
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1581
$ws.Range("J107").Value = 1675.25
$ws.Range("L107").Value = 5025.75
$ws.Range("N107").Value = -8865.75
$ws.Range("H122").Value = 4310.1562
$ws.Range("I122").Value = 819.75
$ws.Range("J122").Value = 4808.7856
$ws.Range("K122").Value = 7377.75
$ws.Range("L122").Value = 43279.0704
$ws.Range("M122").Value = -4927.75
$ws.Range("N122").Value = -48179.0704
$ws.Range("H132").Value = 26375.35
$ws.Range("I132").Value = 1231.25
$ws.Range("J132").Value = 32661.375
$ws.Range("K132").Value = 11081.25
$ws.Range("L132").Value = 293952.375
$ws.Range("M132").Value = -8551.25
$ws.Range("N132").Value = -299012.375

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4049.6511
$ws.Range("I132").Value = 4187.1943
$ws.Range("K132").Value = 12561.5829
$ws.Range("M132").Value = -10031.5829
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
$ws.Range("H139").Value = 59666.668
$ws.Range("J139").Value = 59666.668
$ws.Range("L139").Value = 59666.668
$ws.Range("N139").Value = -69946.66800000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4446.0527
$ws.Range("I122").Value = 4317
$ws.Range("K122").Value = 12951
$ws.Range("M122").Value = -10501
$ws.Range("H136").Value = 3921.9666
$ws.Range("I136").Value = 2455.5625
$ws.Range("J136").Value = 5597.857
$ws.Range("K136").Value = 7366.6875
$ws.Range("L136").Value = 16793.571
$ws.Range("M136").Value = -4816.6875
$ws.Range("N136").Value = -21893.571

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 21551.732
$ws.Range("I107").Value = 1768.9231
$ws.Range("K107").Value = 5306.7693
$ws.Range("M107").Value = -3386.7693
$ws.Range("H122").Value = 4027.923
$ws.Range("I122").Value = 2120.2285
$ws.Range("K122").Value = 6360.685500000001
$ws.Range("M122").Value = -3910.685500000001
$ws.Range("H132").Value = 16793.482
$ws.Range("I132").Value = 24160.53
$ws.Range("K132").Value = 72481.59
$ws.Range("M132").Value = -69951.59
$ws.Range("H139").Value = 138999.5
$ws.Range("J139").Value = 67999
$ws.Range("L139").Value = 67999
$ws.Range("N139").Value = -78279
